# Commit: "Updated fitting parameters and ready tu run detail tests."
#
# The two fitted parameters (r_s_star in column J, h_p_star in column K)
# on the "Parameters" sheet, row 2, were refit and need their new values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parameters")

$ws.Range("J2").Value = 0.04073
$ws.Range("K2").Value = 0.02525
